# "final updation of the product"
# Append a new customer row (customer3 / ford / kljjkl / lkjjkl / kljjkl / lkjjkl)
# to the bottom of the customer table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = "customer3"
$ws.Range("B4").Value = "ford"
$ws.Range("C4").Value = "kljjkl"
$ws.Range("D4").Value = "lkjjkl"
$ws.Range("E4").Value = "kljjkl"
$ws.Range("F4").Value = "lkjjkl"
